$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 7 (FAPs->Neutrophils and all Neutrophils-> rows)
$ws.Range("A4:T7").EntireRow.Delete() | Out-Null

# Update row 2 values (FAPs, Ccl20, Ackr4, ECs)
$ws.Range("I2").Value2 = 1
$ws.Range("J2").Value2 = 1
$ws.Range("M2").Value2 = 0.1631145
$ws.Range("N2").Value2 = 0.326229
$ws.Range("O2").Value2 = 0.7212828052797984
$ws.Range("P2").Value2 = 0.7212828052797984
$ws.Range("Q2").Value2 = 0.1662484188885
$ws.Range("R2").Value2 = 0.9974905133309999
$ws.Range("S2").Value2 = 0.7212828052797984
$ws.Range("T2").Value2 = 0.7212828052797984

# Update row 3 values (FAPs, Ccl20, Ackr4, MuSCs)
$ws.Range("I3").Value2 = 1
$ws.Range("J3").Value2 = 1
$ws.Range("M3").Value2 = 0.0630305
$ws.Range("N3").Value2 = 0.126061
$ws.Range("O3").Value2 = 0.2787171947202017
$ws.Range("P3").Value2 = 0.2787171947202017
$ws.Range("Q3").Value2 = 0.0642415049965
$ws.Range("R3").Value2 = 0.385449029979
$ws.Range("S3").Value2 = 0.2787171947202017
$ws.Range("T3").Value2 = 0.2787171947202017
